# Append a new statistics row to the "統計" (Statistics) sheet,
# mirroring the existing rows' layout: timestamp, count, type,
# skill_match_rate, no_skill_match, multi_skill_match, high_priority.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("統計")

# Determine the next empty row right after the current used range.
$usedRange = $ws.UsedRange
$newRow = $usedRange.Row + $usedRange.Rows.Count

$ws.Cells.Item($newRow, 1).Value = "2025-08-29T12:45:02.484907"
$ws.Cells.Item($newRow, 2).Value = 9
$ws.Cells.Item($newRow, 3).Value = "全案件リスト"
$ws.Cells.Item($newRow, 4).Value = 66.7
$ws.Cells.Item($newRow, 5).Value = 3
$ws.Cells.Item($newRow, 6).Value = 3
$ws.Cells.Item($newRow, 7).Value = 9
